$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Bump the "datetimeFigureOut" date field text on the slide master and all
#    slide layouts from 3/5/2019 -> 3/8/2019.
# ---------------------------------------------------------------------------
function Set-DatePlaceholderText {
    param($container, $text)
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $text
        }
    }
}

Set-DatePlaceholderText $p.SlideMaster "3/8/2019"
for ($l = 1; $l -le $p.SlideMaster.CustomLayouts.Count; $l++) {
    Set-DatePlaceholderText $p.SlideMaster.CustomLayouts.Item($l) "3/8/2019"
}

# ---------------------------------------------------------------------------
# 2. Re-layout / re-colour the Earth cut-away diagram on slide 1.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# Shape index 1 (id=4)
$sh = $s.Shapes.Item(1)
$sh.Left = 5.536732213385827
$sh.Top = 6.255472470866142
$sh.Width = 527.4891662582677
$sh.Height = 527.4891662582677
$sh.VerticalFlip = -1
$sh.Fill.ForeColor.RGB = 0x336699

# Shape index 2 (id=5)
$sh = $s.Shapes.Item(2)
$sh.Left = 7.541456722834646
$sh.Top = 8.26019668031496
$sh.Width = 524.1398926196849
$sh.Height = 524.1398926196849
$sh.VerticalFlip = -1

# Shape index 3 (id=6)
$sh = $s.Shapes.Item(3)
$sh.Left = 23.494291368503937
$sh.Top = 24.54310996614173
$sh.Width = 491.5740509480315
$sh.Height = 491.5740509480315

# Shape index 4 (id=7)
$sh = $s.Shapes.Item(4)
$sh.Left = 38.768306716535434
$sh.Top = 39.817125314173225
$sh.Width = 461.0260315519685
$sh.Height = 461.0260315519685
$sh.VerticalFlip = -1

# Shape index 5 (id=12)
$sh = $s.Shapes.Item(5)
$sh.Left = 317.2679901559055
$sh.Top = 184.9512177023622
$sh.Width = 107.06830601653544
$sh.Height = 107.06830601653544

# Shape index 6 (id=10)
$sh = $s.Shapes.Item(6)
$sh.Left = 277.5779876559055
$sh.Top = 193.09027102047244
$sh.Width = 157.43751521496063
$sh.Height = 157.43751521496063

# Shape index 7 (id=8)
$sh = $s.Shapes.Item(7)
$sh.Left = 121.86334608661417
$sh.Top = 122.91216658425198
$sh.Width = 294.835937511811
$sh.Height = 294.835937511811

# Shape index 8 (id=9)
$sh = $s.Shapes.Item(8)
$sh.Left = 218.252014203937
$sh.Top = 219.3008270015748
$sh.Width = 102.05862047716535
$sh.Height = 102.05862047716535

# Shape index 9 (id=19)
$sh = $s.Shapes.Item(9)
$sh.Left = 531.874145588189
$sh.Top = 255.45941157874017
$sh.Width = 117.02476120944883
$sh.Height = 29.081298842519686
$sh.TextFrame.TextRange.Text = "Crust"
$sh.TextFrame.TextRange.Font.Color.RGB = 0x336699

# Shape index 10 (id=20)
$sh = $s.Shapes.Item(10)
$sh.Left = 531.874145588189
$sh.Top = 280.2424011047244
$sh.Width = 120.6024818448819
$sh.Height = 29.081298842519686

# Shape index 11 (id=21)
$sh = $s.Shapes.Item(11)
$sh.Left = 531.874145588189
$sh.Top = 302.9423981047244
$sh.Width = 130.40823367637796
$sh.Height = 29.081298842519686

# Shape index 12 (id=22)
$sh = $s.Shapes.Item(12)
$sh.Left = 531.874145588189
$sh.Top = 327.3257140913386
$sh.Width = 119.84177018346458
$sh.Height = 29.081298842519686

# Shape index 13 (id=23)
$sh = $s.Shapes.Item(13)
$sh.Left = 531.6366577732283
$sh.Top = 351.18791201574805
$sh.Width = 59.68381882755906
$sh.Height = 29.081298842519686

# Shape index 14 (id=24)
$sh = $s.Shapes.Item(14)
$sh.Left = 531.2166748732283
$sh.Top = 373.8675995551181
$sh.Width = 52.999094018110235
$sh.Height = 29.081298842519686

# Shape index 15 (id=25)
$sh = $s.Shapes.Item(15)
$sh.Left = 531.2166748732283
$sh.Top = 395.81562805118114
$sh.Width = 129.4084624968504
$sh.Height = 29.081298842519686

# Shape index 16 (id=26)
$sh = $s.Shapes.Item(16)
$sh.Left = 531.2166748732283
$sh.Top = 418.49531559055123
$sh.Width = 93.01917265826772
$sh.Height = 29.081298842519686

